$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string must be forced to
# Text format first, so Excel keeps them as text (matching the source data)
# instead of auto-converting them to numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

$ws.Range("D2").Value = "35.529.10"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "1.911.58"
$ws.Range("E3").Value = "  +2.95%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "247.26"
$ws.Range("E5").Value = "  +4.29%  "
$ws.Range("D6").Value = "0.657"
$ws.Range("E6").Value = "  +5.41%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "42.28"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "0.347"
$ws.Range("E9").Value = "  +5.21%  "
$ws.Range("D10").Value = "48.99"
$ws.Range("E10").Value = "  +5.04%  "
$ws.Range("D11").Value = "0.0719"
$ws.Range("E11").Value = "  +3.40%  "
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "2.189.82"
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("D14").Value = "12.35"
$ws.Range("E14").Value = "  +7.86%  "
$ws.Range("D15").Value = "0.702"
$ws.Range("E15").Value = "  +3.46%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "4.87"
$ws.Range("E16").Value = "  +3.56%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.891.49"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "35.560.85"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "72.47"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("E20").Value = "  +4.67%  "
$ws.Range("D21").Value = "244.69"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("D22").Value = "12.74"
$ws.Range("E22").Value = "  +4.92%  "
$ws.Range("E23").Value = "  +1.91%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  +17.23%  "
$ws.Range("D27").Value = "171.64"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "8.51"
$ws.Range("E28").Value = "  +7.02%  "
$ws.Range("D29").Value = "18.36"
$ws.Range("E29").Value = "  +3.75%  "
$ws.Range("E30").Value = "  +3.78%  "
$ws.Range("D31").Value = "0.976"
$ws.Range("E31").Value = "  +24.65%  "
$ws.Range("E32").Value = "  +4.12%  "
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("D34").Value = "4.21"
$ws.Range("E34").Value = "  +4.90%  "
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").Value = "1.72"
$ws.Range("E36").Value = "  +5.92%  "
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").Value = "1.33"
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("D40").Value = "0.0208"
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("D41").Value = "93.13"
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("D42").Value = "0.0634"
$ws.Range("E42").Value = "  +15.17%  "
$ws.Range("D43").Value = "15.63"
$ws.Range("E43").Value = "  +5.09%  "
$ws.Range("D44").Value = "1.351.95"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "2.41"
$ws.Range("E45").Value = "  +2.72%  "
$ws.Range("D46").Value = "47.80"
$ws.Range("E46").Value = "  +39.91%  "
$ws.Range("D47").Value = "12.73"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("D51").Value = "2.100.32"
$ws.Range("E51").Value = "  +3.02%  "
